$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell formatting (bold, border, centered) from H1 onto the
# new I1/J1 header cells before setting their text, so the new cells reuse
# the existing style rather than creating a brand new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @{
    2  = @(10, 10)
    3  = @(9, 9)
    4  = @(4, 5)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(8, 8)
    8  = @(7, 7)
    9  = @(7, 7)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(5, 5)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(7, 8)
    16 = @(7, 7)
    17 = @(6, 6)
    18 = @(8, 8)
    19 = @(7, 7)
    20 = @(7, 7)
    21 = @(5, 5)
    22 = @(3, 3)
    23 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
